$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 0.1706104473901454
$ws.Range("D2").Value = 62.31053733252715
$ws.Range("E2").Value = 0.1955853598633503
$ws.Range("F2").Value = 3.870118989977402
$ws.Range("C3").Value = 0.2298023984463549
$ws.Range("D3").Value = 59.56397644224442
$ws.Range("E3").Value = 0.5642044283991474
$ws.Range("F3").Value = 4.196436727309679
$ws.Range("C4").Value = 0.2737137586955004
$ws.Range("D4").Value = 56.7040827067207
$ws.Range("E4").Value = 0.799883865053875
$ws.Range("F4").Value = 4.333702635138771
$ws.Range("C5").Value = 0.4399171853926251
$ws.Range("D5").Value = 55.41985076582901
$ws.Range("E5").Value = 0.9432651724641165
$ws.Range("F5").Value = 4.607631160917641
$ws.Range("C6").Value = 0.6483481266198544
$ws.Range("D6").Value = 51.80902694792792
$ws.Range("E6").Value = 0.8850631916060058
$ws.Range("F6").Value = 4.797902504061089
$ws.Range("C7").Value = 0.9406170850212482
$ws.Range("D7").Value = 50.58331222310121
$ws.Range("E7").Value = 0.6803566277856896
$ws.Range("F7").Value = 5.102326693539771
$ws.Range("C8").Value = 1.193427487836374
$ws.Range("D8").Value = 46.70925082713161
$ws.Range("E8").Value = 0.446642398865443
$ws.Range("F8").Value = 5.333963642337987
$ws.Range("C9").Value = 1.501492924420091
$ws.Range("D9").Value = 46.31113117345381
$ws.Range("E9").Value = 0.2615607863627488
$ws.Range("F9").Value = 5.281052955962818
$ws.Range("C10").Value = 1.545675494877139
$ws.Range("D10").Value = 46.11917431191578
$ws.Range("E10").Value = 0.1832869415889947
$ws.Range("F10").Value = 5.852505077893012
$ws.Range("C11").Value = 1.740762375146364
$ws.Range("D11").Value = 41.89316015034417
$ws.Range("E11").Value = 0.1589767977762712
$ws.Range("F11").Value = 5.713683919955233
$ws.Range("C12").Value = 1.869093421330333
$ws.Range("D12").Value = 41.0301781710164
$ws.Range("E12").Value = 0.1344422697076893
$ws.Range("F12").Value = 6.038717277985183
$ws.Range("C13").Value = 1.907932412020461
$ws.Range("D13").Value = 38.81282455958083
$ws.Range("E13").Value = 0.1289425282015715
$ws.Range("F13").Value = 6.506188296898556
$ws.Range("C14").Value = 1.990549595035623
$ws.Range("D14").Value = 38.6347565925314
$ws.Range("E14").Value = 0.1027806766296796
$ws.Range("F14").Value = 6.468445766914844
$ws.Range("C15").Value = 2.162627205915117
$ws.Range("D15").Value = 36.01679590681084
$ws.Range("E15").Value = 0.08213908323139289
$ws.Range("F15").Value = 6.957505137613103
$ws.Range("C16").Value = 2.271315168315922
$ws.Range("D16").Value = 36.66317278488118
$ws.Range("E16").Value = 0.05935998359254825
$ws.Range("F16").Value = 6.917161454671909
$ws.Range("C17").Value = 2.31469846353643
$ws.Range("D17").Value = 36.81839981135299
$ws.Range("E17").Value = 0.1058364207803489
$ws.Range("F17").Value = 7.327268325167203
$ws.Range("C18").Value = 2.406487335589338
$ws.Range("D18").Value = 32.86518231563931
$ws.Range("E18").Value = 0.07322132419431246
$ws.Range("F18").Value = 7.497048815969332
$ws.Range("C19").Value = 2.508249329861055
$ws.Range("D19").Value = 33.82676606016329
$ws.Range("E19").Value = 0.07542436753909443
$ws.Range("F19").Value = 7.515381959960477
$ws.Range("C20").Value = 2.528987828299717
$ws.Range("D20").Value = 30.56173728144743
$ws.Range("E20").Value = 0.06467137130463441
$ws.Range("F20").Value = 7.815976653321488
$ws.Range("C21").Value = 2.610139832498959
$ws.Range("D21").Value = 31.08163920554193
$ws.Range("E21").Value = 0.08149811986726199
$ws.Range("F21").Value = 8.283423324529839
